# Auto-generated Excel COM-interop script
# Applies meteocat data refresh update (2026-02-27 03:50 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-27 03:48:17"
$ws.Range("N2").Value = "0.8 °C 3:20 TU"
$ws.Range("O2").Value = "1.6 °C"
$ws.Range("E3").Value = "2026-02-27 03:48:20"
$ws.Range("H3").Value = "35%"
$ws.Range("N3").Value = "1.8 °C 3:14 TU"
$ws.Range("O3").Value = "3.0 °C"
$ws.Range("E4").Value = "2026-02-27 03:48:22"
$ws.Range("H4").Value = "99%"
$ws.Range("J4").Value = "1026.2 hPa"
$ws.Range("N4").Value = "6.4 °C 3:20 TU"
$ws.Range("O4").Value = "7.1 °C"
$ws.Range("E5").Value = "2026-02-27 03:48:25"
$ws.Range("L5").Value = "8.3 km/h - 230º 3:01 TU"
$ws.Range("N5").Value = "2.5 °C 3:08 TU"
$ws.Range("O5").Value = "3.4 °C"
$ws.Range("E6").Value = "2026-02-27 03:48:27"
$ws.Range("J6").Value = "1026.0 hPa"
$ws.Range("N6").Value = "9.0 °C 3:00 TU"
$ws.Range("E7").Value = "2026-02-27 03:48:29"
$ws.Range("J7").Value = "1026.3 hPa"
$ws.Range("N7").Value = "9.6 °C 3:05 TU"
$ws.Range("O7").Value = "10.5 °C"
$ws.Range("E8").Value = "2026-02-27 03:48:32"
$ws.Range("J8").Value = "1025.8 hPa"
$ws.Range("L8").Value = "25.2 km/h - 256º 3:20 TU"
$ws.Range("M8").Value = "12.4 °C 3:19 TU"
$ws.Range("O8").Value = "11.8 °C"
$ws.Range("E9").Value = "2026-02-27 03:48:34"
$ws.Range("O9").Value = "8.1 °C"
$ws.Range("E10").Value = "2026-02-27 03:48:37"
$ws.Range("E11").Value = "2026-02-27 03:48:39"
$ws.Range("N11").Value = "2.0 °C 3:21 TU"
$ws.Range("O11").Value = "2.8 °C"
$ws.Range("E12").Value = "2026-02-27 03:48:41"
$ws.Range("M12").Value = "8.6 °C 3:29 TU"
$ws.Range("O12").Value = "7.4 °C"
$ws.Range("E13").Value = "2026-02-27 03:48:44"
$ws.Range("J13").Value = "1032.2 hPa"
$ws.Range("E14").Value = "2026-02-27 03:48:46"
$ws.Range("N14").Value = "9.0 °C 3:29 TU"
$ws.Range("O14").Value = "9.5 °C"
$ws.Range("E15").Value = "2026-02-27 03:48:48"
$ws.Range("M15").Value = "9.1 °C 3:29 TU"
$ws.Range("O15").Value = "7.9 °C"
$ws.Range("E16").Value = "2026-02-27 03:48:51"
$ws.Range("E17").Value = "2026-02-27 03:48:53"
$ws.Range("N17").Value = "7.0 °C 3:17 TU"
$ws.Range("E18").Value = "2026-02-27 03:48:55"
$ws.Range("J18").Value = "1026.1 hPa"
$ws.Range("N18").Value = "9.5 °C 3:28 TU"
$ws.Range("E19").Value = "2026-02-27 03:48:58"
$ws.Range("E20").Value = "2026-02-27 03:49:00"
$ws.Range("L20").Value = "23.8 km/h - 277º 3:28 TU"
$ws.Range("O20").Value = "1.4 °C"
$ws.Range("E21").Value = "2026-02-27 03:49:02"
$ws.Range("J21").Value = "1029.2 hPa"
$ws.Range("N21").Value = "2.7 °C 3:26 TU"
$ws.Range("O21").Value = "4.3 °C"
$ws.Range("E22").Value = "2026-02-27 03:49:05"
$ws.Range("H22").Value = "49%"
$ws.Range("N22").Value = "-0.2 °C 3:21 TU"
$ws.Range("O22").Value = "0.9 °C"
$ws.Range("E23").Value = "2026-02-27 03:49:07"
$ws.Range("E24").Value = "2026-02-27 03:49:10"
$ws.Range("N24").Value = "2.8 °C 3:24 TU"
$ws.Range("O24").Value = "5.6 °C"
$ws.Range("E25").Value = "2026-02-27 03:49:12"
$ws.Range("H25").Value = "25%"
$ws.Range("O25").Value = "4.8 °C"
$ws.Range("E26").Value = "2026-02-27 03:49:14"
$ws.Range("J26").Value = "1025.2 hPa"
$ws.Range("E27").Value = "2026-02-27 03:49:17"
$ws.Range("M27").Value = "4.5 °C 3:20 TU"
$ws.Range("N27").Value = "3.7 °C 3:07 TU"
$ws.Range("E28").Value = "2026-02-27 03:49:19"
$ws.Range("J28").Value = "1026.3 hPa"
$ws.Range("L28").Value = "6.5 km/h - 262º 3:05 TU"
$ws.Range("N28").Value = "5.2 °C 3:27 TU"
$ws.Range("O28").Value = "6.0 °C"
$ws.Range("E29").Value = "2026-02-27 03:49:22"
$ws.Range("H29").Value = "98%"
$ws.Range("L29").Value = "7.9 km/h - 171º 3:24 TU"
$ws.Range("M29").Value = "10.4 °C 3:29 TU"
$ws.Range("O29").Value = "9.2 °C"
$ws.Range("E30").Value = "2026-02-27 03:49:24"
$ws.Range("J30").Value = "1025.9 hPa"
$ws.Range("O30").Value = "10.1 °C"
$ws.Range("E31").Value = "2026-02-27 03:49:27"
$ws.Range("J31").Value = "1025.5 hPa"
$ws.Range("E32").Value = "2026-02-27 03:49:29"
$ws.Range("H32").Value = "95%"
$ws.Range("N32").Value = "0.7 °C 3:24 TU"
$ws.Range("O32").Value = "1.4 °C"
$ws.Range("E33").Value = "2026-02-27 03:49:32"
$ws.Range("N33").Value = "1.8 °C 3:25 TU"
$ws.Range("O33").Value = "3.0 °C"
$ws.Range("E34").Value = "2026-02-27 03:49:34"
$ws.Range("H34").Value = "48%"
$ws.Range("O34").Value = "1.6 °C"
$ws.Range("E35").Value = "2026-02-27 03:49:36"
$ws.Range("H35").Value = "43%"
$ws.Range("J35").Value = "1025.7 hPa"
$ws.Range("E36").Value = "2026-02-27 03:49:39"
$ws.Range("J36").Value = "1026.2 hPa"
$ws.Range("M36").Value = "9.9 °C 3:29 TU"
$ws.Range("O36").Value = "8.7 °C"
$ws.Range("E37").Value = "2026-02-27 03:49:41"
$ws.Range("L37").Value = "15.8 km/h - 243º 3:16 TU"
$ws.Range("N37").Value = "2.3 °C 3:25 TU"
$ws.Range("O37").Value = "3.0 °C"
$ws.Range("E38").Value = "2026-02-27 03:49:44"
$ws.Range("L38").Value = "9.0 km/h - 295º 3:29 TU"
$ws.Range("N38").Value = "7.4 °C 3:29 TU"
$ws.Range("O38").Value = "8.0 °C"
$ws.Range("E39").Value = "2026-02-27 03:49:46"
$ws.Range("L39").Value = "23.0 km/h - 302º 3:21 TU"
$ws.Range("M39").Value = "5.9 °C 3:19 TU"
$ws.Range("O39").Value = "5.1 °C"
$ws.Range("E40").Value = "2026-02-27 03:49:48"
$ws.Range("N40").Value = "1.6 °C 3:20 TU"
$ws.Range("O40").Value = "2.4 °C"
$ws.Range("E41").Value = "2026-02-27 03:49:51"
$ws.Range("J41").Value = "1026.3 hPa"
$ws.Range("N41").Value = "7.8 °C 3:29 TU"
$ws.Range("O41").Value = "9.5 °C"
$ws.Range("E42").Value = "2026-02-27 03:49:53"
$ws.Range("M42").Value = "9.0 °C 3:27 TU"
$ws.Range("O42").Value = "7.9 °C"
$ws.Range("E43").Value = "2026-02-27 03:49:55"
$ws.Range("N43").Value = "3.3 °C 3:29 TU"
$ws.Range("O43").Value = "4.6 °C"
$ws.Range("E44").Value = "2026-02-27 03:49:57"
$ws.Range("O44").Value = "-0.5 °C"
$ws.Range("E45").Value = "2026-02-27 03:50:00"
$ws.Range("J45").Value = "1026.7 hPa"
$ws.Range("N45").Value = "5.6 °C 3:29 TU"
$ws.Range("E46").Value = "2026-02-27 03:50:02"
$ws.Range("N46").Value = "5.9 °C 3:21 TU"
$ws.Range("O46").Value = "8.0 °C"
